# Apply the "Add data for 2022-04-18" update to the carjacking-by-neighborhood
# -by-month workbook: rename the sheet / refresh the "through" label for the
# current month, and bump a handful of historical April counts (one new
# incident recorded against each neighborhood/year that had a carjacking
# fall on the newly-added day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet title + running "through" label -------------------------------
$ws.Name = "Through 2022-04-10"
$ws.Range("B1").Value = "April 2022 (through April 10)"

# --- Row 2: Austin ---------------------------------------------------------
$ws.Range("F2").Value = 2
$ws.Range("R2").Value = 4
$ws.Range("AD2").Value = 1

# --- Row 4: North Lawndale -------------------------------------------------
$ws.Range("F4").Value = 1
$ws.Range("V4").Value = 5

# --- Row 5: Garfield Park ---------------------------------------------------
$ws.Range("R5").Value = 3

# --- Row 11: Chatham --------------------------------------------------------
$ws.Range("F11").Value = 1

# --- Row 16: Washington Heights ---------------------------------------------
$ws.Range("J16").Value = 2

# --- Row 17: Belmont Cragin --------------------------------------------------
$ws.Range("B17").Value = 1

# --- Row 18: Woodlawn --------------------------------------------------------
$ws.Range("F18").Value = 1

# --- Row 19: Lake View --------------------------------------------------------
$ws.Range("B19").Value = 1

# --- Row 23: Auburn Gresham ---------------------------------------------------
$ws.Range("B23").Value = 2

# --- Row 26: South Shore -------------------------------------------------------
$ws.Range("F26").Value = 1

# --- Row 27: Uptown -------------------------------------------------------------
$ws.Range("B27").Value = 1

# --- Row 29: West Town -----------------------------------------------------------
$ws.Range("V29").Value = 1

# --- Row 42: Avondale --------------------------------------------------------------
$ws.Range("Z42").Value = 1

# --- Row 44: Brighton Park -----------------------------------------------------------
$ws.Range("Z44").Value = 2

# --- Row 45: Douglas ------------------------------------------------------------------
$ws.Range("B45").Value = 1

# --- Row 47: Gage Park ------------------------------------------------------------------
$ws.Range("AD47").Value = 1

# --- Row 70: Hermosa ----------------------------------------------------------------------
$ws.Range("R70").Value = 1

# --- Row 76: Montclare -----------------------------------------------------------------------
$ws.Range("B76").Value = 2

# --- Row 82: Printers Row --------------------------------------------------------------------
$ws.Range("F82").Value = 2

# --- Row 87: South Chicago --------------------------------------------------------------------
$ws.Range("V87").Value = 1

# --- Row 88: Streeterville ----------------------------------------------------------------------
$ws.Range("V88").Value = 1

# --- Row 92: West Pullman ------------------------------------------------------------------------
$ws.Range("F92").Value = 1
